# Update cryptocurrency price/volume data per upstream GitHub Actions refresh
# (commit: "Updated cryptos list on Thu May 23 09:47:22 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '69.671.33'
$ws.Cells.Item(2, 5).Value = '  -0.39%  '
# Row 3
$ws.Cells.Item(3, 4).Value = '3.810.75'
$ws.Cells.Item(3, 5).Value = '  +1.98%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.17%  '
# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '612.10'
$ws.Cells.Item(5, 5).Value = '  -1.68%  '
# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '176.37'
$ws.Cells.Item(6, 5).Value = '  -1.96%  '
# Row 7
$ws.Cells.Item(7, 4).Value = '3.813.05'
$ws.Cells.Item(7, 5).Value = '  +2.10%  '
# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.10%  '
# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.526'
$ws.Cells.Item(9, 5).Value = '  -1.58%  '
# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.167'
$ws.Cells.Item(10, 5).Value = '  -0.34%  '
# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.46'
$ws.Cells.Item(11, 5).Value = '  +2.32%  '
# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.482'
$ws.Cells.Item(12, 5).Value = '  -0.76%  '
# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '39.64'
$ws.Cells.Item(13, 5).Value = '  -3.18%  '
# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.0000253'
$ws.Cells.Item(14, 5).Value = '  -2.29%  '
# Row 15
$ws.Cells.Item(15, 4).Value = '4.455.32'
$ws.Cells.Item(15, 5).Value = '  +2.12%  '
# Row 16
$ws.Cells.Item(16, 4).Value = '3.825.49'
$ws.Cells.Item(16, 5).Value = '  +2.31%  '
# Row 17
$ws.Cells.Item(17, 4).Value = '69.775.14'
$ws.Cells.Item(17, 5).Value = '  -0.32%  '
# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '7.52'
$ws.Cells.Item(18, 5).Value = '  -0.43%  '
# Row 19
$ws.Cells.Item(19, 5).Value = '  -3.21%  '
# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '16.59'
$ws.Cells.Item(20, 5).Value = '  -1.17%  '
# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '505.53'
$ws.Cells.Item(21, 5).Value = '  -0.02%  '
# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '9.55'
$ws.Cells.Item(22, 5).Value = '  +1.95%  '
# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.739'
$ws.Cells.Item(23, 5).Value = '  +2.26%  '
# Row 24
$ws.Cells.Item(24, 2).Value = 'Litecoin'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '86.14'
$ws.Cells.Item(24, 5).Value = '  -0.52%  '
# Row 25
$ws.Cells.Item(25, 2).Value = 'Fetch.AI'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.45'
$ws.Cells.Item(25, 5).Value = '  -2.68%  '
# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.0000141'
$ws.Cells.Item(26, 5).Value = '  +3.28%  '
# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.61'
$ws.Cells.Item(27, 5).Value = '  -3.51%  '
# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.43'
$ws.Cells.Item(28, 5).Value = '  -6.60%  '
# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.09%  '
# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.52'
$ws.Cells.Item(30, 5).Value = '  +1.96%  '
# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.97'
$ws.Cells.Item(31, 5).Value = '  +1.47%  '
# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.96'
$ws.Cells.Item(32, 5).Value = '  +0.59%  '
# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '31.66'
$ws.Cells.Item(33, 5).Value = '  +1.73%  '
# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.113'
$ws.Cells.Item(34, 5).Value = '  -1.58%  '
# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 5).Value = '  +0.12%  '
# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.04'
$ws.Cells.Item(36, 5).Value = '  -1.43%  '
# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '6.09'
$ws.Cells.Item(37, 5).Value = '  -1.81%  '
# Row 38
$ws.Cells.Item(38, 5).Value = '  +5.42%  '
# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '482.96'
$ws.Cells.Item(39, 5).Value = '  +13.62%  '
# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.335'
$ws.Cells.Item(40, 5).Value = '  +0.14%  '
# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.01'
$ws.Cells.Item(41, 5).Value = '  +6.46%  '
# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.05'
$ws.Cells.Item(42, 5).Value = '  -2.75%  '
# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '49.74'
$ws.Cells.Item(43, 5).Value = '  -1.45%  '
# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '43.67'
$ws.Cells.Item(44, 5).Value = '  -2.95%  '
# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '8.53'
$ws.Cells.Item(45, 5).Value = '  -1.96%  '
# Row 46
$ws.Cells.Item(46, 4).Value = '2.919.89'
$ws.Cells.Item(46, 5).Value = '  -2.52%  '
# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0360'
$ws.Cells.Item(47, 5).Value = '  -1.03%  '
# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '139.47'
$ws.Cells.Item(48, 5).Value = '  +1.13%  '
# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '26.87'
$ws.Cells.Item(50, 5).Value = '  -1.50%  '
# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.42'
$ws.Cells.Item(51, 5).Value = '  -3.40%  '
